# RICK_YR_FIN.xlsx update — add a new most-recent fiscal-year column (FY2018,
# period ending 2018-09-30) to the Income Statement / Balance Sheet / Cash
# Flow Statement tables by inserting a new column D and filling it with the
# new period's figures. All pre-existing year columns shift one column right
# (D->E, E->F, ... J->K, K->L).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at D; everything from D onward shifts right.
$ws.Range("D1").EntireColumn.Insert()

# The engine's column-insert copies formatting from the column to the left
# (C), but we want the new column D to look like the data columns that just
# got pushed into E (date format in the "Period Ending" rows, number format
# elsewhere). Copy formats only from E into D to fix this up.
$ws.Range("E1:E200").Copy()
$ws.Range("D1:D200").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# ---- Income Statement (rows 7-35) ----
$ws.Range("D7").Value2 = 43373          # Period Ending -> 30-Sep-18
$ws.Range("D8").Value2 = 165700         # Total Revenue
$ws.Range("D9").Value2 = 67500          # Cost of Revenue
$ws.Range("D10").Value2 = 98300         # Gross Profit
$ws.Range("D12").Value2 = "NA"          # Research Development
$ws.Range("D13").Value2 = 0             # Selling General and Administrative
$ws.Range("D14").Value2 = 6400          # Non Recurring
$ws.Range("D15").Value2 = 7700          # Others
$ws.Range("D17").Value2 = 137400        # Total Operating Expenses
$ws.Range("D18").Value2 = 28400         # Operating Income or Loss
$ws.Range("D20").Value2 = 200           # Total Other Income/Expenses Net
$ws.Range("D21").Value2 = 36400         # Earnings Before Interest And Taxes
$ws.Range("D22").Value2 = 10000         # Interest Expense
$ws.Range("D23").Value2 = 18700         # Income Before Tax
$ws.Range("D24").Value2 = 5600          # Income Tax Expense
$ws.Range("D25").Value2 = 0             # Minority Interest
$ws.Range("D26").Value2 = 13100         # Income After Tax
$ws.Range("D27").Value2 = 13000         # Net Income From Continuing Ops
$ws.Range("D28").Value2 = 0             # Non-recurring Events
$ws.Range("D29").Value2 = 8700          # Discontinued Operations
$ws.Range("D30").Value2 = 0             # Extraordinary Items
$ws.Range("D31").Value2 = 0             # Effect Of Accounting Changes
$ws.Range("D32").Value2 = -200          # Other Items
$ws.Range("D33").Value2 = 21700         # Net Income
$ws.Range("D34").Value2 = 0             # Preferred Stock And Other Adjustments
$ws.Range("D35").Value2 = 21700         # Net Income Applicable To Common Shares

# ---- Balance Sheet (rows 38-77) ----
$ws.Range("D38").Value2 = 43373         # Period Ending -> 30-Sep-18
$ws.Range("D41").Value2 = 17700         # Cash And Cash Equivalents
$ws.Range("D42").Value2 = "NA"          # Short Term Investments
$ws.Range("D43").Value2 = 7300          # Net Receivables
$ws.Range("D44").Value2 = 2400          # Inventory
$ws.Range("D45").Value2 = 9400          # Other Current Assets
$ws.Range("D46").Value2 = 36800         # Total Current Assets
$ws.Range("D47").Value2 = 2900          # Long Term Investments
$ws.Range("D48").Value2 = 344800        # Property Plant and Equipment
$ws.Range("D49").Value2 = 115100        # Goodwill
$ws.Range("D50").Value2 = 0             # Intangible Assets
$ws.Range("D51").Value2 = 0             # Accumulated Amortization
$ws.Range("D52").Value2 = 2500          # Other Assets
$ws.Range("D53").Value2 = 0             # Deferred Long Term Asset Charges
$ws.Range("D54").Value2 = 329700        # Total Assets
$ws.Range("D57").Value2 = 2800          # Accounts Payable
$ws.Range("D58").Value2 = 19000         # Short/Current Long Term Debt
$ws.Range("D59").Value2 = 12000         # Other Current Liabilities
$ws.Range("D60").Value2 = 33800         # Total Current Liabilities
$ws.Range("D61").Value2 = 121600        # Long Term Debt
$ws.Range("D62").Value2 = 21000         # Other Liabilities
$ws.Range("D63").Value2 = 0             # Deferred Long Term Liability Charges
$ws.Range("D64").Value2 = 0             # Minority Interest
$ws.Range("D65").Value2 = 0             # Negative Goodwill
$ws.Range("D66").Value2 = 176300        # Total Liabilities
$ws.Range("D68").Value2 = 0             # Misc Stocks Options Warrants
$ws.Range("D69").Value2 = 0             # Redeemable Preferred Stock
$ws.Range("D70").Value2 = 0             # Preferred Stock
$ws.Range("D71").Value2 = 0             # Common Stock
$ws.Range("D72").Value2 = 88900         # Retained Earnings
$ws.Range("D73").Value2 = 0             # Treasury Stock
$ws.Range("D74").Value2 = 0             # Capital Surplus
$ws.Range("D75").Value2 = 0             # Other Stockholder Equity
$ws.Range("D76").Value2 = 153400        # Total Stockholder Equity
$ws.Range("D77").Value2 = 0             # Net Tangible Assets

# ---- Cash Flow Statement (rows 80-102) ----
$ws.Range("D80").Value2 = 43373         # Period Ending -> 30-Sep-18
$ws.Range("D81").Value2 = 21700         # Net Income
$ws.Range("D83").Value2 = 7700          # Depreciation
$ws.Range("D84").Value2 = 0             # Adjustments To Net Income
$ws.Range("D85").Value2 = 0             # Changes In Accounts Receivables
$ws.Range("D86").Value2 = 0             # Changes In Liabilities
$ws.Range("D87").Value2 = 0             # Changes In Inventories
$ws.Range("D88").Value2 = 0             # Changes In Other Operating Activities
$ws.Range("D89").Value2 = 25800         # Total Cash Flow From Operating Activities
$ws.Range("D91").Value2 = -25300        # Capital Expenditures
$ws.Range("D92").Value2 = 0             # Investments
$ws.Range("D93").Value2 = 0             # Other Cashflows from Investing Activities
$ws.Range("D94").Value2 = -26300        # Total Cash Flows From Investing Activities
$ws.Range("D96").Value2 = -1200         # Dividends Paid
$ws.Range("D97").Value2 = 0             # Sale Purchase of Stock
$ws.Range("D98").Value2 = 0             # Net Borrowings
$ws.Range("D99").Value2 = 0             # Other Cash Flows from Financing Activities
$ws.Range("D100").Value2 = 8400         # Total Cash Flows From Financing Activities
$ws.Range("D101").Value2 = 0            # Effect Of Exchange Rate Changes
$ws.Range("D102").Value2 = 7800         # Change In Cash and Cash Equivalents
